$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value, taken from the crypto-price refresh diff.
# Values are written with a leading apostrophe so Excel stores them as literal
# text (matching the source workbook's inlineStr cells) instead of re-parsing
# number-looking strings (e.g. "0.200" -> 0.2, "9.20" -> 9.2). ClearFormats()
# afterwards drops the quote-prefix cell style Excel adds for that entry method,
# so the on-disk style stays identical to the untouched cells.
$updates = [ordered]@{
    'D2' = '90.497.08'
    'E2' = '  -0.92%  '
    'D3' = '3.134.11'
    'E3' = '  +0.91%  '
    'D4' = '0.998'
    'E4' = '  -0.18%  '
    'D5' = '214.85'
    'E5' = '  -2.00%  '
    'D6' = '622.63'
    'E6' = '  +0.55%  '
    'D7' = '1.12'
    'E7' = '  +24.09%  '
    'D8' = '0.363'
    'E8' = '  -4.35%  '
    'D9' = '0.999'
    'E9' = '  -0.04%  '
    'D10' = '3.132.13'
    'E10' = '  +0.98%  '
    'D11' = '0.744'
    'E11' = '  +9.02%  '
    'D12' = '0.200'
    'E12' = '  +5.70%  '
    'D13' = '5.65'
    'E13' = '  +4.63%  '
    'D14' = '0.0000243'
    'E14' = '  -5.04%  '
    'D15' = '34.90'
    'E15' = '  +4.91%  '
    'D16' = '90.326.62'
    'E16' = '  -0.71%  '
    'D17' = '3.723.46'
    'E17' = '  +1.35%  '
    'D18' = '3.162.08'
    'E18' = '  +3.48%  '
    'D19' = '3.71'
    'E19' = '  +2.14%  '
    'D20' = '14.52'
    'E20' = '  +4.70%  '
    'D21' = '463.61'
    'E21' = '  +6.95%  '
    'D22' = '0.0000210'
    'E22' = '  -7.32%  '
    'D23' = '9.06'
    'E23' = '  +5.63%  '
    'D24' = '5.34'
    'E24' = '  +3.53%  '
    'D25' = '95.31'
    'E25' = '  +13.40%  '
    'D26' = '5.73'
    'E26' = '  +1.81%  '
    'D27' = '12.20'
    'E27' = '  +2.13%  '
    'D28' = '3.319.17'
    'E28' = '  +1.75%  '
    'D29' = '0.999'
    'E29' = '  +0.08%  '
    'D30' = '0.162'
    'E30' = '  -2.96%  '
    'B31' = 'InternetComputer(DFINITY)'
    'C31' = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
    'D31' = '9.20'
    'E31' = '  +5.26%  '
    'B32' = 'Stellar'
    'C32' = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    'D32' = '0.211'
    'E32' = '  +44.93%  '
    'B33' = 'EthereumClassic'
    'C33' = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
    'D33' = '26.52'
    'E33' = '  +14.86%  '
    'B34' = 'Bittensor'
    'C34' = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
    'D34' = '514.74'
    'E34' = '  -1.03%  '
    'B35' = 'PancakeSwap'
    'C35' = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
    'D35' = '1.93'
    'E35' = '  +5.07%  '
    'B36' = 'Kaspa'
    'C36' = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
    'D36' = '0.145'
    'E36' = '  +2.54%  '
    'B37' = 'RenderToken'
    'C37' = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
    'D37' = '6.96'
    'E37' = '  -1.11%  '
    'B38' = 'Fetch.AI'
    'C38' = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
    'D38' = '1.32'
    'E38' = '  +2.03%  '
    'B39' = 'dogwifhat'
    'C39' = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
    'D39' = '3.58'
    'E39' = '  -8.28%  '
    'B40' = 'Hedera'
    'C40' = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
    'D40' = '0.0903'
    'E40' = '  +25.76%  '
    'B41' = 'WhiteBITCoin'
    'C41' = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
    'D41' = '22.21'
    'E41' = '  -0.43%  '
    'B42' = 'Binance-PegBSC-USD'
    'C42' = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
    'D42' = '0.758'
    'E42' = '  -24.19%  '
    'B43' = 'PolygonEcosystemToken'
    'C43' = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
    'D43' = '0.424'
    'E43' = '  +14.41%  '
    'D44' = '0.999'
    'E44' = '  -0.14%  '
    'D45' = '1.98'
    'E45' = '  +4.99%  '
    'D46' = '0.732'
    'E46' = '  +20.01%  '
    'E47' = '  +0.02%  '
    'D48' = '4.73'
    'E48' = '  +12.15%  '
    'D49' = '150.48'
    'E49' = '  +5.51%  '
    'D50' = '1.36'
    'E50' = '  +8.61%  '
    'D51' = '45.12'
    'E51' = '  +3.22%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $updates[$addr]
    $cell.ClearFormats()
}
